$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(44, 8).Value = 27500  # H44
$ws.Cells.Item(44, 10).Value = 27500  # J44
$ws.Cells.Item(44, 12).Value = 27500  # L44
$ws.Cells.Item(44, 14).Value = -28424  # N44
$ws.Cells.Item(62, 8).Value = 3874.125  # H62
$ws.Cells.Item(62, 9).Value = 2831.6667  # I62
$ws.Cells.Item(62, 11).Value = 2831.6667  # K62
$ws.Cells.Item(62, 13).Value = -2207.6667  # M62
$ws.Cells.Item(65, 8).Value = 3874.125  # H65
$ws.Cells.Item(65, 9).Value = 2831.6667  # I65
$ws.Cells.Item(65, 11).Value = 14158.3335  # K65
$ws.Cells.Item(65, 13).Value = -11038.3335  # M65
$ws.Cells.Item(69, 8).Value = 6750  # H69
$ws.Cells.Item(69, 9).Value = 3000  # I69
$ws.Cells.Item(69, 10).Value = 7285.7144  # J69
$ws.Cells.Item(69, 11).Value = 9000  # K69
$ws.Cells.Item(69, 12).Value = 21857.1432  # L69
$ws.Cells.Item(69, 13).Value = -8126  # M69
$ws.Cells.Item(69, 14).Value = -23605.1432  # N69
$ws.Cells.Item(72, 8).Value = 6750  # H72
$ws.Cells.Item(72, 9).Value = 3000  # I72
$ws.Cells.Item(72, 10).Value = 7285.7144  # J72
$ws.Cells.Item(72, 11).Value = 27000  # K72
$ws.Cells.Item(72, 12).Value = 65571.4296  # L72
$ws.Cells.Item(72, 13).Value = -22632  # M72
$ws.Cells.Item(72, 14).Value = -74307.4296  # N72
$ws.Cells.Item(132, 8).Value = 4573.8125  # H132
$ws.Cells.Item(132, 9).Value = 1514  # I132
$ws.Cells.Item(132, 11).Value = 4542  # K132
$ws.Cells.Item(132, 13).Value = -2012  # M132
$ws.Cells.Item(137, 8).Value = 14929574  # H137
$ws.Cells.Item(137, 9).Value = 502092.6  # I137
$ws.Cells.Item(137, 11).Value = 1506277.8  # K137
$ws.Cells.Item(137, 13).Value = -1503727.8  # M137
$ws.Cells.Item(138, 8).Value = 5435.2383  # H138
$ws.Cells.Item(138, 9).Value = 1421.1364  # I138
$ws.Cells.Item(138, 11).Value = 4263.4092  # K138
$ws.Cells.Item(138, 13).Value = 876.5908  # M138
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 2999.5  # H3
$ws.Cells.Item(3, 9).Value = 2999.5  # I3
$ws.Cells.Item(3, 11).Value = 2999.5  # K3
$ws.Cells.Item(3, 13).Value = -2884.5  # M3
$ws.Cells.Item(17, 8).Value = 0  # H17
$ws.Cells.Item(17, 10).Value = 0  # J17
$ws.Cells.Item(17, 12).Value = 0  # L17
$ws.Cells.Item(17, 14).ClearContents()  # N17
$ws.Cells.Item(61, 8).Value = 4481.273  # H61
$ws.Cells.Item(61, 9).Value = 4481.273  # I61
$ws.Cells.Item(61, 11).Value = 4481.273  # K61
$ws.Cells.Item(61, 13).Value = -4269.273  # M61
$ws.Cells.Item(74, 8).Value = 27780660  # H74
$ws.Cells.Item(74, 9).Value = 31252932  # I74
$ws.Cells.Item(74, 10).Value = 2500  # J74
$ws.Cells.Item(74, 11).Value = 31252932  # K74
$ws.Cells.Item(74, 12).Value = 2500  # L74
$ws.Cells.Item(74, 13).Value = -31252058  # M74
$ws.Cells.Item(74, 14).Value = -4248  # N74
$ws.Cells.Item(77, 8).Value = 27780660  # H77
$ws.Cells.Item(77, 9).Value = 31252932  # I77
$ws.Cells.Item(77, 10).Value = 2500  # J77
$ws.Cells.Item(77, 11).Value = 156264660  # K77
$ws.Cells.Item(77, 12).Value = 12500  # L77
$ws.Cells.Item(77, 13).Value = -156260292  # M77
$ws.Cells.Item(77, 14).Value = -21236  # N77
$ws.Cells.Item(104, 8).Value = 50741.668  # H104
$ws.Cells.Item(104, 10).Value = 50741.668  # J104
$ws.Cells.Item(104, 12).Value = 50741.668  # L104
$ws.Cells.Item(104, 14).Value = -57729.668  # N104
$ws.Cells.Item(122, 8).Value = 3963.3333  # H122
$ws.Cells.Item(122, 9).Value = 3529.923  # I122
$ws.Cells.Item(122, 11).Value = 10589.769  # K122
$ws.Cells.Item(122, 13).Value = -8139.769  # M122
$ws.Cells.Item(132, 8).Value = 2865.3  # H132
$ws.Cells.Item(132, 9).Value = 2865.3  # I132
$ws.Cells.Item(132, 11).Value = 8595.900000000001  # K132
$ws.Cells.Item(132, 13).Value = -6065.900000000001  # M132
$ws.Cells.Item(136, 8).Value = 4481.273  # H136
$ws.Cells.Item(136, 9).Value = 4481.273  # I136
$ws.Cells.Item(136, 11).Value = 13443.819  # K136
$ws.Cells.Item(136, 13).Value = -10893.819  # M136
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 7725.643  # H20
$ws.Cells.Item(20, 9).Value = 6961  # I20
$ws.Cells.Item(20, 11).Value = 6961  # K20
$ws.Cells.Item(20, 13).Value = -6714  # M20
$ws.Cells.Item(22, 8).Value = 778  # H22
$ws.Cells.Item(22, 10).Value = 0  # J22
$ws.Cells.Item(22, 12).Value = 0  # L22
$ws.Cells.Item(22, 14).ClearContents()  # N22
$ws.Cells.Item(105, 8).Value = 3123.8  # H105
$ws.Cells.Item(105, 9).Value = 3032.238  # I105
$ws.Cells.Item(105, 11).Value = 3032.238  # K105
$ws.Cells.Item(105, 13).Value = -1285.238  # M105
$ws.Cells.Item(134, 8).Value = 3291.4167  # H134
$ws.Cells.Item(134, 9).Value = 2649.6  # I134
$ws.Cells.Item(134, 11).Value = 7948.799999999999  # K134
$ws.Cells.Item(134, 13).Value = -5413.799999999999  # M134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 22224858  # H31
$ws.Cells.Item(31, 9).Value = 23811898  # I31
$ws.Cells.Item(31, 11).Value = 23811898  # K31
$ws.Cells.Item(31, 13).Value = -23811603  # M31
$ws.Cells.Item(34, 8).Value = 22224858  # H34
$ws.Cells.Item(34, 9).Value = 23811898  # I34
$ws.Cells.Item(34, 11).Value = 23811898  # K34
$ws.Cells.Item(34, 13).Value = -23811696  # M34
$ws.Cells.Item(58, 8).Value = 2604.0908  # H58
$ws.Cells.Item(58, 9).Value = 2604.0908  # I58
$ws.Cells.Item(58, 11).Value = 2604.0908  # K58
$ws.Cells.Item(58, 13).Value = -2401.0908  # M58
$ws.Cells.Item(62, 8).Value = 65204.5  # H62
$ws.Cells.Item(62, 9).Value = 39298  # I62
$ws.Cells.Item(62, 11).Value = 39298  # K62
$ws.Cells.Item(62, 13).Value = -38674  # M62
$ws.Cells.Item(65, 8).Value = 65204.5  # H65
$ws.Cells.Item(65, 9).Value = 39298  # I65
$ws.Cells.Item(65, 11).Value = 196490  # K65
$ws.Cells.Item(65, 13).Value = -193370  # M65
$ws.Cells.Item(132, 8).Value = 95254060  # H132
$ws.Cells.Item(132, 9).Value = 133343690  # I132
$ws.Cells.Item(132, 11).Value = 400031070  # K132
$ws.Cells.Item(132, 13).Value = -400028540  # M132
$ws.Cells.Item(134, 8).Value = 3224.6667  # H134
$ws.Cells.Item(134, 9).Value = 2154.2307  # I134
$ws.Cells.Item(134, 10).Value = 4489.727  # J134
$ws.Cells.Item(134, 11).Value = 6462.6921  # K134
$ws.Cells.Item(134, 12).Value = 13469.181  # L134
$ws.Cells.Item(134, 13).Value = -3927.6921  # M134
$ws.Cells.Item(134, 14).Value = -18539.181  # N134
$ws.Cells.Item(136, 8).Value = 2604.0908  # H136
$ws.Cells.Item(136, 9).Value = 2604.0908  # I136
$ws.Cells.Item(136, 11).Value = 7812.2724  # K136
$ws.Cells.Item(136, 13).Value = -5262.2724  # M136
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 374  # H86
$ws.Cells.Item(86, 9).Value = 374  # I86
$ws.Cells.Item(86, 11).Value = 1122  # K86
$ws.Cells.Item(86, 13).Value = 64  # M86
$ws.Cells.Item(89, 8).Value = 374  # H89
$ws.Cells.Item(89, 9).Value = 374  # I89
$ws.Cells.Item(89, 11).Value = 3366  # K89
$ws.Cells.Item(89, 13).Value = 2562  # M89
$ws.Cells.Item(116, 8).Value = 4788.1665  # H116
$ws.Cells.Item(116, 9).Value = 6014.5  # I116
$ws.Cells.Item(116, 11).Value = 18043.5  # K116
$ws.Cells.Item(116, 13).Value = -14601.5  # M116
$ws.Cells.Item(136, 8).Value = 6335.0625  # H136
$ws.Cells.Item(136, 9).Value = 1422.3334  # I136
$ws.Cells.Item(136, 10).Value = 9282.7  # J136
$ws.Cells.Item(136, 11).Value = 4267.0002  # K136
$ws.Cells.Item(136, 12).Value = 27848.1  # L136
$ws.Cells.Item(136, 13).Value = 832.9997999999996  # M136
$ws.Cells.Item(136, 14).Value = -38048.10000000001  # N136
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 1851146.4  # H11
$ws.Cells.Item(11, 9).Value = 806780.6  # I11
$ws.Cells.Item(11, 10).Value = 2503875  # J11
$ws.Cells.Item(11, 11).Value = 806780.6  # K11
$ws.Cells.Item(11, 12).Value = 2503875  # L11
$ws.Cells.Item(11, 13).Value = -806641.6  # M11
$ws.Cells.Item(11, 14).Value = -2504153  # N11
$ws.Cells.Item(18, 8).Value = 17000  # H18
$ws.Cells.Item(18, 9).Value = 17000  # I18
$ws.Cells.Item(18, 11).Value = 17000  # K18
$ws.Cells.Item(18, 13).Value = -16707  # M18
$ws.Cells.Item(102, 8).Value = 10003488  # H102
$ws.Cells.Item(102, 9).Value = 11631508  # I102
$ws.Cells.Item(102, 11).Value = 11631508  # K102
$ws.Cells.Item(102, 13).Value = -11629886  # M102
$ws.Cells.Item(126, 8).Value = 2581.6155  # H126
$ws.Cells.Item(126, 9).Value = 2361.2  # I126
$ws.Cells.Item(126, 10).Value = 3316.3333  # J126
$ws.Cells.Item(126, 11).Value = 7083.599999999999  # K126
$ws.Cells.Item(126, 12).Value = 9948.999899999999  # L126
$ws.Cells.Item(126, 13).Value = -4613.599999999999  # M126
$ws.Cells.Item(126, 14).Value = -14888.9999  # N126
$ws.Cells.Item(132, 8).Value = 254341.62  # H132
$ws.Cells.Item(132, 9).Value = 668144.3  # I132
$ws.Cells.Item(132, 11).Value = 2004432.9  # K132
$ws.Cells.Item(132, 13).Value = -2001902.9  # M132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1675.8667  # H93
$ws.Cells.Item(93, 9).Value = 1779.6  # I93
$ws.Cells.Item(93, 10).Value = 1468.4  # J93
$ws.Cells.Item(93, 11).Value = 1779.6  # K93
$ws.Cells.Item(93, 12).Value = 1468.4  # L93
$ws.Cells.Item(93, 13).Value = -531.5999999999999  # M93
$ws.Cells.Item(93, 14).Value = -3964.4  # N93
$ws.Cells.Item(94, 8).Value = 0  # H94
$ws.Cells.Item(94, 10).Value = 0  # J94
$ws.Cells.Item(94, 12).Value = 0  # L94
$ws.Cells.Item(94, 14).ClearContents()  # N94
$ws.Cells.Item(105, 8).Value = 44000  # H105
$ws.Cells.Item(105, 10).Value = 44000  # J105
$ws.Cells.Item(105, 12).Value = 44000  # L105
$ws.Cells.Item(105, 14).Value = -50988  # N105
$ws.Cells.Item(122, 8).Value = 6771.8687  # H122
$ws.Cells.Item(122, 9).Value = 3890.7083  # I122
$ws.Cells.Item(122, 10).Value = 11711  # J122
$ws.Cells.Item(122, 11).Value = 11672.1249  # K122
$ws.Cells.Item(122, 12).Value = 35133  # L122
$ws.Cells.Item(122, 13).Value = -9222.124899999999  # M122
$ws.Cells.Item(122, 14).Value = -40033  # N122
$ws.Cells.Item(136, 8).Value = 3935.4  # H136
$ws.Cells.Item(136, 9).Value = 2732.5264  # I136
$ws.Cells.Item(136, 11).Value = 8197.5792  # K136
$ws.Cells.Item(136, 13).Value = -5647.5792  # M136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 29036.25  # H41
$ws.Cells.Item(41, 10).Value = 29036.25  # J41
$ws.Cells.Item(41, 12).Value = 29036.25  # L41
$ws.Cells.Item(41, 14).Value = -29816.25  # N41
$ws.Cells.Item(81, 8).Value = 9493.8  # H81
$ws.Cells.Item(81, 9).Value = 2487.5  # I81
$ws.Cells.Item(81, 10).Value = 14164.667  # J81
$ws.Cells.Item(81, 11).Value = 4975  # K81
$ws.Cells.Item(81, 12).Value = 28329.334  # L81
$ws.Cells.Item(81, 13).Value = -3914  # M81
$ws.Cells.Item(81, 14).Value = -30451.334  # N81
$ws.Cells.Item(84, 8).Value = 9493.8  # H84
$ws.Cells.Item(84, 9).Value = 2487.5  # I84
$ws.Cells.Item(84, 10).Value = 14164.667  # J84
$ws.Cells.Item(84, 11).Value = 24875  # K84
$ws.Cells.Item(84, 12).Value = 141646.67  # L84
$ws.Cells.Item(84, 13).Value = -19571  # M84
$ws.Cells.Item(84, 14).Value = -152254.67  # N84
$ws.Cells.Item(132, 8).Value = 1495.8243  # H132
$ws.Cells.Item(132, 9).Value = 1523.1818  # I132
$ws.Cells.Item(132, 10).Value = 1491.0476  # J132
$ws.Cells.Item(132, 11).Value = 4569.5454  # K132
$ws.Cells.Item(132, 12).Value = 4473.142800000001  # L132
$ws.Cells.Item(132, 13).Value = -2039.5454  # M132
$ws.Cells.Item(132, 14).Value = -9533.142800000001  # N132
